$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Input")

# New Order IDs replacing the previous test values in column R (OrderId)
# for rows 2, 3, 5 and 6. The cells are text-typed (they hold numeric-
# looking order id strings), so write through a Text-formatted helper
# cell and PasteSpecial the *values* across — this preserves the
# original "General" number format / style of each target cell while
# still landing the content as text rather than letting Excel coerce a
# numeric-looking string into a Number.

$helper = $ws.Range("Z1")

$helper.NumberFormat = "@"
$helper.Value = "51484912"
$helper.Copy()
$ws.Range("R2").PasteSpecial(-4163)

$helper.NumberFormat = "@"
$helper.Value = "51484913"
$helper.Copy()
$ws.Range("R3").PasteSpecial(-4163)

$helper.NumberFormat = "@"
$helper.Value = "51484882"
$helper.Copy()
$ws.Range("R5").PasteSpecial(-4163)

$helper.NumberFormat = "@"
$helper.Value = "51484915"
$helper.Copy()
$ws.Range("R6").PasteSpecial(-4163)

$helper.Clear()
$excel.CutCopyMode = $false
